$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric (e.g. "1.007", "29.083.49").
# Force the Price column to Text format first so Excel stores the updated
# values as literal text instead of auto-converting them to numbers,
# matching the original inline-string formatting used for this column.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.083.49"
$ws.Range("E2").Value = "  -0.39%  "

# Row 3
$ws.Range("D3").Value = "1.815.04"
$ws.Range("E3").Value = "  -0.80%  "

# Row 4
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.68%  "

# Row 5
$ws.Range("D5").Value = "232.71"
$ws.Range("E5").Value = "  -2.13%  "

# Row 6
$ws.Range("D6").Value = "0.5834"
$ws.Range("E6").Value = "  -3.61%  "

# Row 7
$ws.Range("D7").Value = "1.007"
$ws.Range("E7").Value = "  +0.60%  "

# Row 8
$ws.Range("D8").Value = "0.2695"
$ws.Range("E8").Value = "  -4.74%  "

# Row 9
$ws.Range("D9").Value = "0.06700"
$ws.Range("E9").Value = "  -5.54%  "

# Row 10
$ws.Range("D10").Value = "22.61"
$ws.Range("E10").Value = "  -5.75%  "

# Row 11
$ws.Range("D11").Value = "0.07508"
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("D12").Value = "1.811.43"
$ws.Range("E12").Value = "  -0.88%  "

# Row 13
$ws.Range("D13").Value = "4.591"
$ws.Range("E13").Value = "  -4.36%  "

# Row 14
$ws.Range("D14").Value = "0.6107"
$ws.Range("E14").Value = "  -4.26%  "

# Row 15
$ws.Range("D15").Value = "0.000009324"
$ws.Range("E15").Value = "  -6.42%  "

# Row 16
$ws.Range("D16").Value = "73.94"
$ws.Range("E16").Value = "  -7.28%  "

# Row 17
$ws.Range("D17").Value = "28.862.88"
$ws.Range("E17").Value = "  -0.96%  "

# Row 18
$ws.Range("D18").Value = "5.365"
$ws.Range("E18").Value = "  -10.41%  "

# Row 19
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.57%  "

# Row 20
$ws.Range("D20").Value = "205.02"
$ws.Range("E20").Value = "  -11.04%  "

# Row 21
$ws.Range("D21").Value = "11.23"
$ws.Range("E21").Value = "  -4.67%  "

# Row 22
$ws.Range("D22").Value = "6.672"
$ws.Range("E22").Value = "  -4.51%  "

# Row 23
$ws.Range("D23").Value = "1.010"
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
$ws.Range("D24").Value = "154.09"
$ws.Range("E24").Value = "  -0.97%  "

# Row 25
$ws.Range("D25").Value = "7.684"
$ws.Range("E25").Value = "  -4.43%  "

# Row 26
$ws.Range("D26").Value = "0.1246"
$ws.Range("E26").Value = "  -3.37%  "

# Row 27
$ws.Range("D27").Value = "16.05"
$ws.Range("E27").Value = "  -3.93%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "1.406"
$ws.Range("E28").Value = "  -3.31%  "

# Row 29
$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").Value = "0.06259"
$ws.Range("E29").Value = "  -6.23%  "

# Row 30
$ws.Range("D30").Value = "1.433"

# Row 31
$ws.Range("D31").Value = "3.664"
$ws.Range("E31").Value = "  -3.89%  "

# Row 32
$ws.Range("D32").Value = "3.620"
$ws.Range("E32").Value = "  -5.58%  "

# Row 33
$ws.Range("D33").Value = "1.675"
$ws.Range("E33").Value = "  -2.53%  "

# Row 34
$ws.Range("D34").Value = "1.041"
$ws.Range("E34").Value = "  -8.25%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.535"
$ws.Range("E35").Value = "  -0.54%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.6244"
$ws.Range("E36").Value = "  -4.92%  "

# Row 37
$ws.Range("D37").Value = "2.746"
$ws.Range("E37").Value = "  -0.36%  "

# Row 38
$ws.Range("D38").Value = "0.01695"
$ws.Range("E38").Value = "  -3.98%  "

# Row 39
$ws.Range("D39").Value = "6.347"
$ws.Range("E39").Value = "  -3.43%  "

# Row 40
$ws.Range("D40").Value = "1.119.85"
$ws.Range("E40").Value = "  -9.64%  "

# Row 41
$ws.Range("D41").Value = "0.8590"
$ws.Range("E41").Value = "  -7.47%  "

# Row 42
$ws.Range("D42").Value = "1.007"
$ws.Range("E42").Value = "  +0.59%  "

# Row 43
$ws.Range("D43").Value = "1.964.72"

# Row 44
$ws.Range("D44").Value = "99.35"
$ws.Range("E44").Value = "  -0.96%  "

# Row 45
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  -3.54%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "59.63"
$ws.Range("E46").Value = "  -6.00%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.4543"
$ws.Range("E47").Value = "  -0.45%  "

# Row 48
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.558"
$ws.Range("E48").Value = "  -4.40%  "

# Row 49
$ws.Range("D49").Value = "0.05493"
$ws.Range("E49").Value = "  -1.62%  "

# Row 50
$ws.Range("D50").Value = "8.207"
$ws.Range("E50").Value = "  -3.37%  "

# Row 51
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  +0.55%  "
